# Updated UI variable mapping for IAM arrays in Physical trough model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Widen column F to fit the new, longer descriptions
$ws.Columns.Item(6).ColumnWidth = 51

# Insert 15 new rows below the existing table, copying the formatting
# (styles) of the last existing row (14) down into each new row.
for ($i = 0; $i -lt 15; $i++) {
    $ws.Rows.Item(14).Copy()
    $ws.Rows.Item(15).Insert()
}

# New rows of data (15-29)
$data = @(
    @("Changed name","number","csp_dtr_sca_iam0_1","IAMs_1[0]","Physical Trough Collector Type 1","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_1","IAMs_1[1]","Physical Trough Collector Type 1","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_1","IAMs_1[2]","Physical Trough Collector Type 1","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam0_2","IAMs_2[0]","Physical Trough Collector Type 2","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_2","IAMs_2[1]","Physical Trough Collector Type 2","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_2","IAMs_2[2]","Physical Trough Collector Type 2","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam0_3","IAMs_3[0]","Physical Trough Collector Type 3","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_3","IAMs_3[1]","Physical Trough Collector Type 3","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_3","IAMs_3[2]","Physical Trough Collector Type 3","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam0_4","IAMs_4[0]","Physical Trough Collector Type 4","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam1_4","IAMs_4[1]","Physical Trough Collector Type 4","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","number","csp_dtr_sca_iam2_4","IAMs_4[2]","Physical Trough Collector Type 4","Allowing table/array of IAM coefficients as inputs","Ty      If IAM beceoms available in User Library, may want to make values array instead of set at 3 variables"),
    @("Changed name","array","IamF0","IAM_matrix","Physical Trough Collector Header","combining collector IAM coef. Arrays into 1 output matrix","Ty"),
    @("Changed name","array","IamF1","IAM_matrix","Physical Trough Collector Header","combining collector IAM coef. Arrays into 1 output matrix","Ty"),
    @("Changed name","array","IamF2","IAM_matrix","Physical Trough Collector Header","combining collector IAM coef. Arrays into 1 output matrix","Ty")
)

$startRow = 15
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}

# Extend the "Type" list validation down to row 58
$ws.Range("A2:A58").Validation.Delete()
$ws.Range("A2:A58").Validation.Add(3, 1, 1, "=Types")
$ws.Range("A2:A58").Validation.InCellDropdown = $true
$ws.Range("A2:A58").Validation.ShowInput = $true
$ws.Range("A2:A58").Validation.ShowError = $true

# Move the instructional rounded-rectangle callout shape down to sit next to the new rows
$shp = $ws.Shapes.Item(1)
$shp.Left = 941.25
$shp.Top = 315.75
$shp.Width = 314.25
$shp.Height = 124.5

# Update the active selection to reflect where the user left off editing
$ws.Range("D29").Select()
